# Microgrid_cmd.xlsx — "added PV disconnect to prevent it tanking the MG during faults"
#
# P_req (sheet 1): the requested real-power setpoint series was changed from a
#   flat -1,400,000 W (rows 2-38) / 0 W (rows 39-52) profile to a flat
#   +500,000 W profile, with a single one-off spike to 1,000,000 W at t=11
#   (row 13) — modelling the PV disconnecting (and the MG commanding a
#   positive, rather than negative, real-power reference) when a fault trips
#   the island.
# Q_req (sheet 2): no value changes, just a stale cursor-position update.
# try_island (sheet 3): the islanding flag trace now starts high (t=0,1) and
#   is low for the rest of the run (t=17-50 flipped from 1 to 0), and the
#   series is extended out to t=60 (10 new zero rows) to cover the longer
#   disconnect test.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: P_req
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

$ws1.Range("B2:B38").Value = 500000
$ws1.Range("B39:B52").Value = 500000
$ws1.Range("B13").Value = 1000000

$ws1.PageSetup.Orientation = 1

$ws1.Range("B13").Select()
$excel.ActiveWindow.Zoom = 100

# ---------------------------------------------------------------------
# Sheet 2: Q_req
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("L18").Select()

# ---------------------------------------------------------------------
# Sheet 3: try_island
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()

$ws3.Range("B2").Value = 1
$ws3.Range("B3").Value = 1
$ws3.Range("B19:B52").Value = 0

for ($i = 0; $i -lt 10; $i++) {
    $row = 53 + $i
    $ws3.Cells.Item($row, 1).Value = 51 + $i
    $ws3.Cells.Item($row, 2).Value = 0
}

$ws3.Range("E4:F5").Select()
$ws3.Range("F5").Activate()

Write-Output "done"
